$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.352.00"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "2.064.36"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'234.04"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  +1.28%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'56.81"
$ws.Range("E8").Value = "  -1.10%  "

$ws.Range("D9").Value = "'0.381"
$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("E10").Value = "  +0.24%  "

$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("D12").Value = "2.368.96"
$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("D13").Value = "'14.58"
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").Value = "'20.64"
$ws.Range("E14").Value = "  -2.07%  "

$ws.Range("D15").Value = "'0.776"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").Value = "'5.13"
$ws.Range("E16").Value = "  -2.07%  "

$ws.Range("D17").Value = "2.066.31"
$ws.Range("E17").Value = "  -0.82%  "

$ws.Range("D18").Value = "37.272.42"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").Value = "'6.24"
$ws.Range("E19").Value = "  +4.21%  "

$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").Value = "'226.30"
$ws.Range("E22").Value = "  +1.19%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("E25").Value = "  -1.26%  "

$ws.Range("D26").Value = "'167.35"
$ws.Range("E26").Value = "  +2.84%  "

$ws.Range("D27").Value = "'8.75"
$ws.Range("E27").Value = "  -1.20%  "

$ws.Range("E28").Value = "  +3.64%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'19.08"
$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.127"
$ws.Range("E30").Value = "  -3.49%  "

$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("D32").Value = "'4.45"
$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("D33").Value = "'0.0615"
$ws.Range("E33").Value = "  -1.35%  "

$ws.Range("E34").Value = "  +3.78%  "

$ws.Range("D35").Value = "'2.48"
$ws.Range("E35").Value = "  -2.97%  "

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").Value = "'3.23"
$ws.Range("E38").Value = "  -2.74%  "

$ws.Range("D39").Value = "'5.65"
$ws.Range("E39").Value = "  -4.68%  "

$ws.Range("E40").Value = "  -0.20%  "

$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "'0.0940"
$ws.Range("E41").Value = "  -2.25%  "

$ws.Range("D42").Value = "1.463.91"
$ws.Range("E42").Value = "  -0.82%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'96.01"
$ws.Range("E43").Value = "  +1.13%  "

$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "'4.34"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").Value = "'1.17"
$ws.Range("E45").Value = "  +4.02%  "

$ws.Range("E46").Value = "  +1.08%  "

$ws.Range("E47").Value = "  -1.16%  "

$ws.Range("D48").Value = "'15.02"
$ws.Range("E48").Value = "  -8.56%  "

$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.96"
$ws.Range("E49").Value = "  +0.77%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'7.14"
$ws.Range("E50").Value = "  -2.66%  "

$ws.Range("D51").Value = "2.258.73"
$ws.Range("E51").Value = "  -0.46%  "
